$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.740272164344788
$ws.Range("B1").Value = 3.311477661132812
$ws.Range("C1").Value = 4.733545303344727
$ws.Range("D1").Value = 3.052782773971558
$ws.Range("E1").Value = 1.703563332557678
